$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# ---------------------------------------------------------------------------
# 1. Remove the three unused, empty default sheets (Tabelle3/Tabelle2/Tabelle1)
# ---------------------------------------------------------------------------
$wb.Worksheets("Tabelle3").Delete()
$wb.Worksheets("Tabelle2").Delete()
$wb.Worksheets("Tabelle1").Delete()

# ---------------------------------------------------------------------------
# 2. Monsters sheet - new boss monster + overlay group entry
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets("Monsters")
$ws.Range("A6").Value = 61
$ws.Range("B6").Value = "Untoter Lord"
$ws.Range("C6").Value = "Boss in Manyeyes' castle"

$ws.Range("H8").Value = 93
$ws.Range("I8").Value = "1x Untoter Lord, 1x Untoter Krieger, 3x Untoter Magier"

$ws.Activate()
$ws.Range("B7").Select()

# ---------------------------------------------------------------------------
# 3. New Object Graphics sheet - new undead sprite row + note
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets("New Object Graphics")
$ws.Range("A6").Value = 374
$ws.Range("B6").Value = "Undead (that works with sky, pal4)"
$ws.Range("C6").Value = "Manyeyes' castle 2"
$ws.Range("D6").Value = "2Object3D.amb"
$ws.Range("E6").Value = "New"
$ws.Range("I1").Value = "Also added two overlays 92 and 93 for the door in Manyeyes' castle 2"
# Columns B and C were widened to fit their new (longer) content - "bestFit" autosize.
$ws.Columns("B").ColumnWidth = 31.16666666666667
$ws.Columns("C").ColumnWidth = 16.736979166666668

$ws.Activate()
$ws.Range("H18").Select()

# ---------------------------------------------------------------------------
# 4. Items sheet - two new reward items
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets("Items")
$ws.Range("A16").Value = 417
$ws.Range("B16").Value = "Schwert der Ahnen / Ancestral Sword"
$ws.Range("C16").Value = "Weapon"
$ws.Range("D16").Value = "Reward from first manyeyes' castle boss"

$ws.Range("A17").Value = 418
$ws.Range("B17").Value = "Rüstung der Ahnen / Ancestral Armour"
$ws.Range("C17").Value = "Armor"
$ws.Range("D17").Value = "Reward from first manyeyes' castle boss"

$ws.Activate()
$ws.Range("D21").Select()

# ---------------------------------------------------------------------------
# 5. GlobalVars sheet - new global var entry
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets("GlobalVars")
$ws.Range("A30").Value = "247: Upper boss in Manyeyes' castle killed"

$ws.Activate()
$ws.Range("A31").Select()

# ---------------------------------------------------------------------------
# 6. Todo sheet - fix existing entry text + add new todo item
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets("Todo")
$ws.Range("A8").Value = "Fix 462 event icons (after opening the door, it is a closed chest …)"
$ws.Range("A9").Value = "Change small lizard sprite to not include the sky color!"

$ws.Activate()
$ws.Range("F17").Select()

# ---------------------------------------------------------------------------
# 7. Chests sheet - two new rows for Manyeyes' cellar (463)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets("Chests")
$ws.Range("A22").Value = 152
$ws.Range("B22").Value = "Manyeyes'c cellar (463)"
$ws.Range("C22").Value = "100 Gold, 6 Rationen"

$ws.Range("A23").Value = 153
$ws.Range("B23").Value = "Manyeyes'c cellar (463)"
$ws.Range("C23").Value = "1x Strength Potion, 1x Intelligence Potion, 2x Antidot, 5x Healing Potion IV, 3x Spell Potion III"

# Chests ends up as the active sheet/tab, with C23 selected.
$ws.Activate()
$ws.Range("C23").Select()
